$d = $word.ActiveDocument

# The edit removes the whole paragraph that contains nothing but the
# italicised "Zekaria" run sitting directly under the "ZEC" Heading2
# paragraph (it is a short duplicate/stub of the real "Zekaria" Heading2
# further down, which is NOT italicised and must be left untouched).
#
# Identify it robustly: scan paragraphs back-to-front looking for one
# whose trimmed text is exactly "Zekaria" and whose run formatting is
# italic (the real heading's "Zekaria" is not italic), then delete the
# whole paragraph (not just the text) so the surrounding paragraphs join
# together exactly as in the target diff.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.Trim()
    if ($text -eq "Zekaria" -and $p.Range.Font.Italic -eq -1) {
        $p.Range.Delete()
    }
}
